$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 303, shifting existing row 303..392 down to 304..393
$ws.Rows.Item(303).Insert()

# Populate the newly inserted row 303 with the new weekly data record
$ws.Cells.Item(303, 1).Value = 5
$ws.Cells.Item(303, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(303, 3).Value = "Maule"
$ws.Cells.Item(303, 4).Value = "10/30/2023"
$ws.Cells.Item(303, 5).Value = 7
$ws.Cells.Item(303, 6).Value = 100112017
$ws.Cells.Item(303, 7).Value = "Apio"
$ws.Cells.Item(303, 8).Value = "Americana (o)"
$ws.Cells.Item(303, 9).Value = "Primera"
$ws.Cells.Item(303, 10).Value = 400
$ws.Cells.Item(303, 11).Value = 7000
$ws.Cells.Item(303, 12).Value = 7000
$ws.Cells.Item(303, 13).Value = 7000
$ws.Cells.Item(303, 14).Value = "`$/docena de matas"
$ws.Cells.Item(303, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(303, 16).Value = 1167
$ws.Cells.Item(303, 17).Value = 6
$ws.Cells.Item(303, 18).Value = "Hortaliza"
